$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -3.377448218355161
$ws.Range("C2").Value = 2.15496247269007
$ws.Range("D2").Value = 7.434495881238234

$ws.Range("B3").Value = 0.2448820992536849
$ws.Range("C3").Value = -0.1721910504140167
$ws.Range("D3").Value = -1.008236814482644

$ws.Range("B4").Value = 5.624925926329238
$ws.Range("C4").Value = 0.6459332894185987
$ws.Range("D4").Value = 6.267294698459325

$ws.Range("B5").Value = 6.16826137189288
$ws.Range("C5").Value = -6.935153294182605
$ws.Range("D5").Value = 9.755477202390939

$ws.Range("B6").Value = -1.494625744833378
$ws.Range("C6").Value = -6.157851153201799
$ws.Range("D6").Value = 8.009006717074318

$ws.Range("B7").Value = -0.2554344193826941
$ws.Range("C7").Value = -4.932397146504464
$ws.Range("D7").Value = 2.809173804671983

$ws.Range("B8").Value = -0.8961894651313584
$ws.Range("C8").Value = -4.264270422584227
$ws.Range("D8").Value = 0.1319837197746532

$ws.Range("B9").Value = 4.683070112298715
$ws.Range("C9").Value = -1.149109367878443
$ws.Range("D9").Value = 11.26013612946959

$ws.Range("B10").Value = -10.69158489251918
$ws.Range("C10").Value = -5.600889391489416
$ws.Range("D10").Value = -5.900781225340602

$ws.Range("B11").Value = -6.373529693431146
$ws.Range("C11").Value = 9.024371965595002
$ws.Range("D11").Value = -7.350727472305019

$ws.Range("B12").Value = 0.04395139504043133
$ws.Range("C12").Value = 7.481601158193651
$ws.Range("D12").Value = -10.82795411725803

$ws.Range("B13").Value = -2.445450002465022
$ws.Range("C13").Value = 2.990192558263849
$ws.Range("D13").Value = -3.269703462328233
